$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("新增物品")

# Final desired content for rows 6-10 (column A sequence numbers stay 5..9, untouched).
# "crushed_raw_nickel" (old row 6) and "andesite_alloy_block" (old row 10) are removed from
# the list, so every row from 6 to 10 now holds the data that used to sit one row below it,
# and the two renamed English labels are applied along the way.

# Row 6: was "crushed_raw_nickel" -> becomes "raw_nickel_block" (English name updated).
$ws.Range("B6").Value = "raw_nickel_block"
$ws.Range("C6").Value = "粗镍块"
$ws.Range("D6").Value = "Block Of Raw Nickel"

# Row 7: was "raw_nickel_block" -> becomes "nickel_block" (English name updated).
$ws.Range("B7").Value = "nickel_block"
$ws.Range("C7").Value = "镍块"
$ws.Range("D7").Value = "Block Of Nickel"

# Row 8: was "nickel_block" -> becomes "andesite_alloy_sheet".
$ws.Range("B8").Value = "andesite_alloy_sheet"
$ws.Range("C8").Value = "安山合金板"
$ws.Range("D8").Value = "Andesite Alloy Sheet"

# Row 9: was "andesite_alloy_sheet" -> becomes "nickel_rich_laterite".
$ws.Range("B9").Value = "nickel_rich_laterite"
$ws.Range("C9").Value = "富镍红土"
$ws.Range("D9").Value = "Nickel-rich Laterite"

# Row 10: was "andesite_alloy_block" -> becomes "laterite".
$ws.Range("B10").Value = "laterite"
$ws.Range("C10").Value = "红土"
$ws.Range("D10").Value = "Laterite"

# Remove the now-surplus trailing rows 11 and 12.
$ws.Range("A11:A12").EntireRow.Delete()

# Update the view: zoom and new selection.
$win = $excel.ActiveWindow
$win.Zoom = 115
$ws.Range("C14").Select()

$wb.Save()
